# Generate Report for Handback
# The 0822cacf-... file has been handed back (in sync with en-US), so its
# row moves from "Ready for handoff" to "Handed back" status, gains a
# Latest Handback DateTime, and its status flips to "Include". Since the
# report lists files in handback order, 0822cacf now sorts ahead of
# a19ccab5-... (rows 2 and 3 swap across all three sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.md"
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"

$overview.Range("A3").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.md"
$zhcn.Range("B2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.zh-cn.xlf"
$zhcn.Range("D2").Value = "2016-03-08 12:34:08"
$zhcn.Range("E2").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.md"
$zhcn.Range("F2").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.zh-cn.xlf"
$zhcn.Range("G2").Value = "2016-03-08 12:34:31"
$zhcn.Range("H2").Value = "Include"

$zhcn.Range("A3").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md"
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.77927fa78aecdc12cbc3d27452998e4801193aa1.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-08 12:32:42"
$zhcn.Range("E3").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md"
$zhcn.Range("F3").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.77927fa78aecdc12cbc3d27452998e4801193aa1.zh-cn.xlf"
$zhcn.Range("G3").Value = "2016-03-08 12:33:24"
$zhcn.Range("H3").Value = "Include"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.md"
$dede.Range("B2").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.de-de.xlf"
$dede.Range("D2").Value = "2016-03-08 12:34:12"
$dede.Range("E2").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.md"
$dede.Range("F2").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.de-de.xlf"
$dede.Range("G2").Value = "2016-03-08 12:34:41"
$dede.Range("H2").Value = "Include"

$dede.Range("A3").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md"
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.77927fa78aecdc12cbc3d27452998e4801193aa1.de-de.xlf"
$dede.Range("D3").Value = "2016-03-08 12:32:55"
$dede.Range("E3").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md"
$dede.Range("F3").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.77927fa78aecdc12cbc3d27452998e4801193aa1.de-de.xlf"
$dede.Range("G3").Value = "2016-03-08 12:33:35"
$dede.Range("H3").Value = "Include"
